$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 12: finish filling out the existing "10" row with the new Submit Comment test case
$ws.Range("B12").Value = "test_<SubmitComment>"
$ws.Range("C12").Value = "This is to test whether users are able to submit a comment"
$ws.Range("D12").Value = "Name = ""John Low"" `nComment = ""Your Resume is interesting"""
$ws.Range("E12").Value = "Comment is displayed under the 'Comments' secton in the blog page"
$ws.Range("F12").Value = "Comment is displayed under the 'Comments' secton in the blog page"
$ws.Range("G12").Value = "Based on the given source code, this function has already been implemented"

# Wrap text on D12 and bump row height for the multi-line justification text
$ws.Range("D12").WrapText = $true
$ws.Rows.Item(12).RowHeight = 29.15

# Row 13: new Submit Empty Comment test case
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "test_<SubmitEmptyComment>"
$ws.Range("C13").Value = "This is to test whether users are able to submit an empty comment"
$ws.Range("D13").Value = "NIL"
$ws.Range("E13").Value = "Error displayed. Require users to fill in the 2 textboxes"
$ws.Range("F13").Value = "Error displayed. Require users to fill in the 2 textboxes"
$ws.Range("G13").Value = "Not able to test this function in pytest because error is only displayed on client side"

# Row 14: new blank row with just the index counter
$ws.Range("A14").Value = 12

# Update the selection to match the author's final cursor position
$ws.Range("F17").Select()
